$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 4390.4
$ws.Range("I64").Value = 4664
$ws.Range("J64").Value = 3980
$ws.Range("K64").Value = 4664
$ws.Range("L64").Value = 3980
$ws.Range("M64").Value = -4416
$ws.Range("N64").Value = -4476

# Row 67
$ws.Range("H67").Value = 4390.4
$ws.Range("I67").Value = 4664
$ws.Range("J67").Value = 3980
$ws.Range("K67").Value = 4664
$ws.Range("L67").Value = 3980
$ws.Range("M67").Value = -3806
$ws.Range("N67").Value = -5696

# Row 92
$ws.Range("H92").Value = 1580.1666
$ws.Range("I92").Value = 1637.1538
$ws.Range("J92").Value = 1432
$ws.Range("K92").Value = 1637.1538
$ws.Range("L92").Value = 1432
$ws.Range("M92").Value = -389.1538
$ws.Range("N92").Value = -3928

# Row 116
$ws.Range("H116").Value = 2332.8462
$ws.Range("I116").Value = 2158.4707
$ws.Range("J116").Value = 2662.2222
$ws.Range("K116").Value = 2158.4707
$ws.Range("L116").Value = 2662.2222
$ws.Range("M116").Value = 1283.5293
$ws.Range("N116").Value = -9546.2222

# Row 132
$ws.Range("H132").Value = 8780284
$ws.Range("I132").Value = 13895431
$ws.Range("J132").Value = 11460.571
$ws.Range("K132").Value = 41686293
$ws.Range("L132").Value = 34381.713
$ws.Range("M132").Value = -41683763
$ws.Range("N132").Value = -39441.713

# Row 133
$ws.Range("H133").Value = 34942.25
$ws.Range("J133").Value = 34942.25
$ws.Range("L133").Value = 34942.25
$ws.Range("N133").Value = -45062.25

# Row 136
$ws.Range("H136").Value = 34880
$ws.Range("J136").Value = 34880
$ws.Range("L136").Value = 34880
$ws.Range("N136").Value = -45080

# Row 137
$ws.Range("H137").Value = 1647.5264
$ws.Range("I137").Value = 1014.8
$ws.Range("K137").Value = 3044.4
$ws.Range("M137").Value = -494.3999999999996

# Row 138
$ws.Range("H138").Value = 1510.6837
$ws.Range("I138").Value = 818.85
$ws.Range("J138").Value = 1688.0769
$ws.Range("K138").Value = 2456.55
$ws.Range("L138").Value = 5064.2307
$ws.Range("M138").Value = 2683.45
$ws.Range("N138").Value = -15344.2307

# Row 141
$ws.Range("H141").Value = 655.9091
$ws.Range("I141").Value = 655.9091
$ws.Range("K141").Value = 1967.7273
$ws.Range("M141").Value = 3212.2727

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4893.216
$ws.Range("I32").Value = 4370.4565
$ws.Range("K32").Value = 4370.4565
$ws.Range("M32").Value = -4083.4565

# Row 43
$ws.Range("H43").Value = 7822.5
$ws.Range("J43").Value = 7822.5
$ws.Range("L43").Value = 7822.5
$ws.Range("N43").Value = -8448.5

# Row 74
$ws.Range("H74").Value = 1752.3
$ws.Range("I74").Value = 917.875
$ws.Range("J74").Value = 2705.9285
$ws.Range("K74").Value = 917.875
$ws.Range("L74").Value = 2705.9285
$ws.Range("M74").Value = -43.875
$ws.Range("N74").Value = -4453.9285

# Row 77
$ws.Range("H77").Value = 1752.3
$ws.Range("I77").Value = 917.875
$ws.Range("J77").Value = 2705.9285
$ws.Range("K77").Value = 4589.375
$ws.Range("L77").Value = 13529.6425
$ws.Range("M77").Value = -221.375
$ws.Range("N77").Value = -22265.6425

# Row 88
$ws.Range("H88").Value = 2431
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 2586.1667
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 2586.1667
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -3398.1667

# Row 91
$ws.Range("H91").Value = 2431
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 2586.1667
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 2586.1667
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -5394.1667

# Row 122
$ws.Range("H122").Value = 1128.1904
$ws.Range("I122").Value = 934.25
$ws.Range("K122").Value = 2802.75
$ws.Range("M122").Value = -352.75

# Row 138
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

# Row 139
$ws.Range("H139").Value = 49476.668
$ws.Range("J139").Value = 49476.668
$ws.Range("L139").Value = 49476.668
$ws.Range("N139").Value = -59756.668

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 23
$ws.Range("H23").Value = 3014
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 3014
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 3014
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -3580

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1539.2941
$ws.Range("I31").Value = 1242.8
$ws.Range("K31").Value = 1242.8
$ws.Range("M31").Value = -947.8

# Row 34
$ws.Range("H34").Value = 1539.2941
$ws.Range("I34").Value = 1242.8
$ws.Range("K34").Value = 1242.8
$ws.Range("M34").Value = -1040.8

# Row 62
$ws.Range("H62").Value = 10528779
$ws.Range("I62").Value = 2458.8235
$ws.Range("K62").Value = 2458.8235
$ws.Range("M62").Value = -1834.8235

# Row 65
$ws.Range("H65").Value = 10528779
$ws.Range("I65").Value = 2458.8235
$ws.Range("K65").Value = 12294.1175
$ws.Range("M65").Value = -9174.1175

# Row 109
$ws.Range("H109").Value = 23362.625
$ws.Range("J109").Value = 23362.625
$ws.Range("L109").Value = 23362.625
$ws.Range("N109").Value = -25442.625

# Row 132
$ws.Range("H132").Value = 1482.1177
$ws.Range("I132").Value = 1136.9756
$ws.Range("K132").Value = 3410.9268
$ws.Range("M132").Value = -880.9268000000002

# Row 134
$ws.Range("H134").Value = 999.1111
$ws.Range("I134").Value = 852.6923
$ws.Range("J134").Value = 1379.8
$ws.Range("K134").Value = 2558.0769
$ws.Range("L134").Value = 4139.4
$ws.Range("M134").Value = -23.07690000000002
$ws.Range("N134").Value = -9209.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 2528.5334
$ws.Range("I81").Value = 1502.6
$ws.Range("J81").Value = 3041.5
$ws.Range("K81").Value = 4507.799999999999
$ws.Range("L81").Value = 9124.5
$ws.Range("M81").Value = -3384.799999999999
$ws.Range("N81").Value = -11370.5

# Row 84
$ws.Range("H84").Value = 2528.5334
$ws.Range("I84").Value = 1502.6
$ws.Range("J84").Value = 3041.5
$ws.Range("K84").Value = 13523.4
$ws.Range("L84").Value = 27373.5
$ws.Range("M84").Value = -7907.4
$ws.Range("N84").Value = -38605.5

# Row 121
$ws.Range("H121").Value = 633.3333
$ws.Range("I121").Value = 450
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 1350
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = -40
$ws.Range("N121").Value = -5620

# Row 131
$ws.Range("H131").Value = 21279634
$ws.Range("J131").Value = 3895.6667
$ws.Range("L131").Value = 11687.0001
$ws.Range("N131").Value = -21767.0001

# Row 136
$ws.Range("H136").Value = 1422.6666
$ws.Range("J136").Value = 4116.5
$ws.Range("L136").Value = 12349.5
$ws.Range("N136").Value = -22549.5

# Row 139
$ws.Range("H139").Value = 1893.7742
$ws.Range("J139").Value = 1699.1538
$ws.Range("L139").Value = 5097.4614
$ws.Range("N139").Value = -15377.4614

# Row 141
$ws.Range("H141").Value = 4166.6665

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 6400
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 6400
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 6400
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -6746

# Row 30
$ws.Range("H30").Value = 6400
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 6400
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 6400
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -6610

# Row 104
$ws.Range("H104").Value = 48500
$ws.Range("J104").Value = 48500
$ws.Range("L104").Value = 48500
$ws.Range("N104").Value = -55488

# Row 107
$ws.Range("H107").Value = 900.17645
$ws.Range("I107").Value = 900.3333
$ws.Range("J107").Value = 899.8
$ws.Range("K107").Value = 900.3333
$ws.Range("L107").Value = 899.8
$ws.Range("M107").Value = 1019.6667
$ws.Range("N107").Value = -4739.8

# Row 141
$ws.Range("H141").Value = 78000
$ws.Range("J141").Value = 78000
$ws.Range("L141").Value = 78000
$ws.Range("N141").Value = -88360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6128.5
$ws.Range("I40").Value = 3222
$ws.Range("J40").Value = 9035
$ws.Range("K40").Value = 3222
$ws.Range("L40").Value = 9035
$ws.Range("M40").Value = -3086
$ws.Range("N40").Value = -9307

# Row 61
$ws.Range("H61").Value = 1695.8334
$ws.Range("I61").Value = 1106.6666
$ws.Range("J61").Value = 2285
$ws.Range("K61").Value = 1106.6666
$ws.Range("L61").Value = 2285
$ws.Range("M61").Value = -904.6666
$ws.Range("N61").Value = -2689

# Row 113
$ws.Range("H113").Value = 1695.8334
$ws.Range("I113").Value = 1106.6666
$ws.Range("J113").Value = 2285
$ws.Range("K113").Value = 1106.6666
$ws.Range("L113").Value = 2285
$ws.Range("M113").Value = 1063.3334
$ws.Range("N113").Value = -6625

# Row 132
$ws.Range("H132").Value = 23479.826
$ws.Range("I132").Value = 1120.4348
$ws.Range("J132").Value = 45839.22
$ws.Range("K132").Value = 3361.3044
$ws.Range("L132").Value = 137517.66
$ws.Range("M132").Value = -831.3044
$ws.Range("N132").Value = -142577.66

# Row 134
$ws.Range("H134").Value = 29371.428
$ws.Range("J134").Value = 29371.428
$ws.Range("L134").Value = 29371.428
$ws.Range("N134").Value = -39511.428

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 17
$ws.Range("H17").Value = 229.75
$ws.Range("J17").Value = 229.75
$ws.Range("L17").Value = 229.75
$ws.Range("N17").Value = -573.75

# Row 100
$ws.Range("H100").Value = 928.5454999999999
$ws.Range("I100").Value = 1195.8334
$ws.Range("K100").Value = 2391.6668
$ws.Range("M100").Value = -1850.6668

# Row 132
$ws.Range("H132").Value = 2339.1875
$ws.Range("I132").Value = 2226.7585
$ws.Range("J132").Value = 3426
$ws.Range("K132").Value = 6680.2755
$ws.Range("L132").Value = 10278
$ws.Range("M132").Value = -4150.2755
$ws.Range("N132").Value = -15338

# Row 136
$ws.Range("H136").Value = 544.1177
$ws.Range("I136").Value = 239.21428
$ws.Range("J136").Value = 1967
$ws.Range("K136").Value = 717.64284
$ws.Range("L136").Value = 5901
$ws.Range("M136").Value = 1832.35716
$ws.Range("N136").Value = -11001
